{"js": "const body = context.document.body;\n\n// Helper: find the first occurrence of `findText` in the document body and\n// replace it with `replaceText`, preserving the surrounding formatting by\n// reusing the found range.\nasync function replaceFirst(findText, replaceText) {\n  const results = body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n  results.items[0].insertText(replaceText, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"Empirical findings offer a complicated picture. Studies provide\n//    evidence that equalization...\" ->\n//    \"Empirical findings offer support for both hypotheses. Equalization...\"\nawait replaceFirst(\n  \"offer a complicated picture. Studies provide evidence that equalization\",\n  \"offer support for both hypotheses. Equalization\"\n);\n\n// 2) \"...depending on the outcome, as well as...\" ->\n//    \"...depending on the outcome considered, as well as...\"\nawait replaceFirst(\n  \"depending on the outcome, as well as\",\n  \"depending on the outcome considered, as well as\"\n);\n\n// 3) \"...as a moderator, as those with higher levels of interest...\" ->\n//    \"...as a moderator, in that those with higher levels of interest...\"\nawait replaceFirst(\n  \"as a moderator, as those with higher\",\n  \"as a moderator, in that those with higher\"\n);\n\n// 4) \"...where some groups are clearly left in information landscapes...\" ->\n//    \"...where some groups are left in information landscapes...\"\nawait replaceFirst(\n  \"some groups are clearly left\",\n  \"some groups are left\"\n);\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $findText,    # FindText\n        $false,       # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        2             # Replace (wdReplaceOne)\n    ) | Out-Null\n}\n\n# 1) \"Empirical findings offer a complicated picture. Studies provide\n#    evidence that equalization...\" ->\n#    \"Empirical findings offer support for both hypotheses. Equalization...\"\nReplace-Text \"offer a complicated picture. Studies provide evidence that equalization\" \"offer support for both hypotheses. Equalization\"\n\n# 2) \"...depending on the outcome, as well as...\" ->\n#    \"...depending on the outcome considered, as well as...\"\nReplace-Text \"depending on the outcome, as well as\" \"depending on the outcome considered, as well as\"\n\n# 3) \"...as a moderator, as those with higher levels of interest...\" ->\n#    \"...as a moderator, in that those with higher levels of interest...\"\nReplace-Text \"as a moderator, as those with higher\" \"as a moderator, in that those with higher\"\n\n# 4) \"...where some groups are clearly left in information landscapes...\" ->\n#    \"...where some groups are left in information landscapes...\"\nReplace-Text \"some groups are clearly left\" \"some groups are left\"\n"}
